# 600 retune 100Hz loop filter
$wb = $excel.ActiveWorkbook

$wsLoop = $wb.Worksheets.Item("LoopFilter LPF")
$wsBranch = $wb.Worksheets.Item("Branch LPF")

# Retune the loop filter cutoff frequency from 50Hz to 100Hz
$wsLoop.Range("B2").Value = 100

# Restore the default selection on "Branch LPF" before switching away from it
$wsBranch.Range("B19").Select()

# Make "LoopFilter LPF" the active sheet, with C13 selected
$wsLoop.Select()
$wsLoop.Range("C13").Select()
